$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest
# snapshot values. Both columns hold plain text in the source sheet
# (e.g. "42.563.99", "  -2.26%  "), so the Price cells are written
# via .Formula with a leading apostrophe -- the same as typing
# '42.563.99 into the cell -- so Excel stores them as literal text
# instead of silently reinterpreting them as numbers (which would
# mangle thousand-separator-less values and drop trailing zeros).

$ws.Range("D2").Formula = "'42.563.99"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Formula = "'2.278.65"
$ws.Range("E3").Value = "  -4.10%  "
$ws.Range("D5").Formula = "'299.99"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("D6").Formula = "'97.36"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("D10").Formula = "'33.66"
$ws.Range("E10").Value = "  -6.38%  "
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Formula = "'50.53"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Formula = "'6.62"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").Formula = "'2.630.61"
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Formula = "'2.279.81"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").Formula = "'42.435.64"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Formula = "'11.35"
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Formula = "'235.17"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("E25").Value = "  -6.03%  "
$ws.Range("D26").Formula = "'2.47"
$ws.Range("E26").Value = "  -5.23%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Formula = "'24.31"
$ws.Range("E28").Value = "  -6.16%  "
$ws.Range("D29").Formula = "'2.18"
$ws.Range("E29").Value = "  -6.40%  "
$ws.Range("D30").Formula = "'163.98"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").Formula = "'33.51"
$ws.Range("E31").Value = "  -8.70%  "
$ws.Range("D32").Formula = "'9.07"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("D35").Formula = "'2.41"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").Formula = "'0.0695"
$ws.Range("E36").Value = "  -5.58%  "
$ws.Range("E37").Value = "  -7.66%  "
$ws.Range("E38").Value = "  -9.17%  "
$ws.Range("D39").Formula = "'16.03"
$ws.Range("E39").Value = "  -12.28%  "
$ws.Range("E40").Value = "  -5.89%  "
$ws.Range("E41").Value = "  -9.06%  "
$ws.Range("D42").Formula = "'0.110"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("D43").Formula = "'2.40"
$ws.Range("E43").Value = "  -8.72%  "
$ws.Range("D44").Formula = "'1.951.88"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Formula = "'0.0280"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Formula = "'17.63"
$ws.Range("E46").Value = "  -10.85%  "
$ws.Range("D47").Formula = "'9.64"
$ws.Range("E47").Value = "  -8.87%  "
$ws.Range("D48").Formula = "'2.81"
$ws.Range("E48").Value = "  -9.77%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").Formula = "'2.500.38"
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("E51").Value = "  -1.78%  "
